# Applies the Coinranking crypto-price refresh for the GitHub Actions run
# completed on Wed Oct 11 16:32:00 UTC 2023.
#
# Every figure in the sheet (price, 1h volume change %) is stored as plain
# text, so each write goes through Set-CellText below. That helper adds a
# leading single quote whenever the new value looks like a plain number
# (e.g. "22.00" or "0.488"), which stops Excel's automatic type detection
# from turning it into a numeric value and keeps it as text exactly as
# scraped from the site.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        [string]$Cell,
        [string]$Text
    )
    $range = $ws.Range($Cell)
    if ($Text -match '^[+-]?\d+(\.\d+)?$') {
        # Looks like a plain number ("22.00", "0.488", ...) - force text
        # so Excel doesn't silently reinterpret/round it as a number.
        $range.Value = "'" + $Text
    } else {
        $range.Value = $Text
    }
}

Set-CellText "D2" '26.797.42'
Set-CellText "E2" '  -2.55%  '

Set-CellText "D3" '1.564.15'
Set-CellText "E3" '  -0.31%  '

Set-CellText "E4" '  +0.32%  '

Set-CellText "D5" '206.18'
Set-CellText "E5" '  -0.78%  '

Set-CellText "D6" '0.488'
Set-CellText "E6" '  -1.86%  '

Set-CellText "E7" '  +0.30%  '

Set-CellText "D8" '22.00'
Set-CellText "E8" '  -0.13%  '

Set-CellText "E9" '  -0.22%  '

Set-CellText "D10" '0.0583'
Set-CellText "E10" '  -1.09%  '

Set-CellText "D11" '0.0866'
Set-CellText "E11" '  -0.08%  '

Set-CellText "D12" '1.786.60'
Set-CellText "E12" '  -0.40%  '

Set-CellText "D13" '1.566.92'
Set-CellText "E13" '  +0.13%  '

Set-CellText "E14" '  -1.92%  '

Set-CellText "E15" '  -1.03%  '

# Row 16/17: Litecoin and WrappedBTC swapped ranking positions
Set-CellText "B16" 'Litecoin'
Set-CellText "C16" 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-CellText "D16" '61.71'
Set-CellText "E16" '  -2.54%  '

Set-CellText "B17" 'WrappedBTC'
Set-CellText "C17" 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-CellText "D17" '26.842.78'
Set-CellText "E17" '  -2.21%  '

Set-CellText "D18" '215.25'
Set-CellText "E18" '  +0.59%  '

Set-CellText "D19" '7.33'
Set-CellText "E19" '  +1.01%  '

Set-CellText "E20" '  -1.50%  '

Set-CellText "E21" '  +0.30%  '

Set-CellText "E22" '  -0.84%  '

Set-CellText "E23" '  -2.32%  '

Set-CellText "E24" '  -0.12%  '

Set-CellText "D25" '152.02'
Set-CellText "E25" '  -1.17%  '

Set-CellText "D26" '6.74'
Set-CellText "E26" '  -1.40%  '

Set-CellText "D27" '14.87'
Set-CellText "E27" '  -1.13%  '

Set-CellText "E28" '  +0.32%  '

Set-CellText "E29" '  -1.43%  '

# Row 30/31: Hedera and PancakeSwap swapped ranking positions
Set-CellText "B30" 'Hedera'
Set-CellText "C30" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-CellText "D30" '0.0462'
Set-CellText "E30" '  -1.70%  '

Set-CellText "B31" 'PancakeSwap'
Set-CellText "C31" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-CellText "D31" '1.11'
Set-CellText "E31" '  -4.04%  '

Set-CellText "D32" '3.16'
Set-CellText "E32" '  -1.35%  '

Set-CellText "D33" '1.386.79'
Set-CellText "E33" '  +1.72%  '

Set-CellText "D34" '2.91'
Set-CellText "E34" '  -1.43%  '

Set-CellText "D35" '1.55'
Set-CellText "E35" '  +0.66%  '

Set-CellText "D36" '2.30'
Set-CellText "E36" '  -0.19%  '

Set-CellText "D37" '0.940'
Set-CellText "E37" '  -3.43%  '

Set-CellText "D38" '0.0162'
Set-CellText "E38" '  -3.06%  '

Set-CellText "E39" '  -1.28%  '

Set-CellText "E40" '  -3.91%  '

Set-CellText "E41" '  +0.31%  '

Set-CellText "E42" '  +3.07%  '

Set-CellText "D43" '5.44'
Set-CellText "E43" '  +2.95%  '

Set-CellText "D44" '1.78'
Set-CellText "E44" '  -0.46%  '

Set-CellText "D45" '2.18'
Set-CellText "E45" '  +1.68%  '

Set-CellText "D46" '63.34'
Set-CellText "E46" '  -1.18%  '

Set-CellText "D47" '1.699.31'

Set-CellText "D48" '85.72'
Set-CellText "E48" '  +0.44%  '

Set-CellText "D49" '0.0₇0968'
Set-CellText "E49" '  +3.61%  '

Set-CellText "D50" '0.0495'
Set-CellText "E50" '  -0.09%  '

Set-CellText "E51" '  -0.75%  '
